$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the header/column "index" -> "i" (this also auto-renames the
# structured table column "testdata[index]" -> "testdata[i]" and re-syncs
# the shared-strings table).
$ws.Range("A1").Value = "i"

# The underlying data column used to be a 1-based running counter
# (1, 2, 3, ...). It becomes 0-based (0, 1, 2, ...), i.e. every value in
# column A (rows 2-503) is decremented by one.
for ($r = 2; $r -le 503; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 - 1
}

# Column A visually narrows now that "index"/3-digit numbers became
# "i"/up-to-3-digit numbers.
$ws.Columns.Item(1).ColumnWidth = 3.166666666666667
